# Auto-generated: update crypto price (D) and volume-change (E) columns
# Values are prefixed with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr cells) instead of auto-
# converting number-like strings (e.g. "1.001") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.004.17"
$ws.Range("E2").Value = "'  +0.31%  "
$ws.Range("D3").Value = "'1.640.97"
$ws.Range("E3").Value = "'  -0.24%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.47%  "
$ws.Range("D5").Value = "'215.14"
$ws.Range("E5").Value = "'  -0.17%  "
$ws.Range("D6").Value = "'0.5097"
$ws.Range("E6").Value = "'  +0.14%  "
$ws.Range("E7").Value = "'  -0.43%  "
$ws.Range("D8").Value = "'0.2585"
$ws.Range("E8").Value = "'  +0.55%  "
$ws.Range("D9").Value = "'0.06362"
$ws.Range("E9").Value = "'  -0.67%  "
$ws.Range("D10").Value = "'19.92"
$ws.Range("E10").Value = "'  +1.59%  "
$ws.Range("D11").Value = "'0.07754"
$ws.Range("E11").Value = "'  -0.30%  "
$ws.Range("D12").Value = "'4.300"
$ws.Range("E12").Value = "'  -0.24%  "
$ws.Range("D13").Value = "'1.635.93"
$ws.Range("E13").Value = "'  -0.75%  "
$ws.Range("D14").Value = "'0.5479"
$ws.Range("D15").Value = "'0.0₅7754"
$ws.Range("E15").Value = "'  -1.44%  "
$ws.Range("D16").Value = "'64.40"
$ws.Range("E16").Value = "'  -0.41%  "
$ws.Range("D17").Value = "'26.018.38"
$ws.Range("E17").Value = "'  +0.05%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "'  -0.41%  "
$ws.Range("D19").Value = "'196.96"
$ws.Range("E19").Value = "'  -0.76%  "
$ws.Range("D20").Value = "'4.466"
$ws.Range("E20").Value = "'  +0.34%  "
$ws.Range("D21").Value = "'9.972"
$ws.Range("E21").Value = "'  -0.12%  "
$ws.Range("D22").Value = "'6.139"
$ws.Range("E22").Value = "'  +1.52%  "
$ws.Range("E23").Value = "'  -0.63%  "
$ws.Range("D24").Value = "'1.891"
$ws.Range("E24").Value = "'  +0.58%  "
$ws.Range("D25").Value = "'142.89"
$ws.Range("D26").Value = "'0.1266"
$ws.Range("E26").Value = "'  +10.29%  "
$ws.Range("D27").Value = "'6.878"
$ws.Range("E27").Value = "'  -0.54%  "
$ws.Range("D28").Value = "'15.62"
$ws.Range("E28").Value = "'  -0.80%  "
$ws.Range("E29").Value = "'  -0.02%  "
$ws.Range("D30").Value = "'0.04894"
$ws.Range("E30").Value = "'  -2.59%  "
$ws.Range("D31").Value = "'3.287"
$ws.Range("E31").Value = "'  +0.63%  "
$ws.Range("D32").Value = "'3.217"
$ws.Range("E32").Value = "'  +0.74%  "
$ws.Range("D33").Value = "'1.557"
$ws.Range("E33").Value = "'  +1.00%  "
$ws.Range("D34").Value = "'2.376"
$ws.Range("E34").Value = "'  +0.50%  "
$ws.Range("D35").Value = "'0.9198"
$ws.Range("E35").Value = "'  +2.74%  "
$ws.Range("D36").Value = "'2.566"
$ws.Range("E36").Value = "'  -0.96%  "
$ws.Range("D37").Value = "'0.5560"
$ws.Range("E37").Value = "'  +0.68%  "
$ws.Range("D38").Value = "'1.105.56"
$ws.Range("E38").Value = "'  -2.52%  "
$ws.Range("D39").Value = "'0.01570"
$ws.Range("E39").Value = "'  +0.49%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "'  -0.52%  "
$ws.Range("D41").Value = "'5.616"
$ws.Range("E41").Value = "'  -0.35%  "
$ws.Range("D42").Value = "'0.8045"
$ws.Range("E42").Value = "'  -1.57%  "
$ws.Range("D43").Value = "'98.76"
$ws.Range("E43").Value = "'  -1.05%  "
$ws.Range("E44").Value = "'  -4.81%  "
$ws.Range("D45").Value = "'1.781.07"
$ws.Range("E45").Value = "'  -0.09%  "
$ws.Range("D46").Value = "'0.4536"
$ws.Range("E46").Value = "'  +0.07%  "
$ws.Range("D47").Value = "'55.40"
$ws.Range("E47").Value = "'  +0.74%  "
$ws.Range("E48").Value = "'  -0.25%  "
$ws.Range("D49").Value = "'0.05186"
$ws.Range("D50").Value = "'7.562"
$ws.Range("E50").Value = "'  +2.07%  "
$ws.Range("E51").Value = "'  -0.23%  "
